$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting of the last existing entry row (row 15) down to the new row (16)
$ws.Range("A15:B15").Copy()
$ws.Range("A16:B16").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Fill in the new log entry values
$ws.Range("A16").Value = "2/22, 3 hours"
$ws.Range("B16").Value = "Worked on shiny app, adding histogram and side layouts"

# Match the row height used by the other wrapped-text log rows (e.g. row 15)
$ws.Rows.Item(16).RowHeight = $ws.Rows.Item(15).RowHeight

# Update selection to reflect where the cursor moved after data entry (B17)
$ws.Range("B17").Select()
